$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Section A6:D9 - turn the single "quantity" column into a per-city
#     load table (Madrid / Barcelona / Lisboa), and move the destination /
#     distance mini-table from D:E to F:G ---

$ws.Range("A6").Value = "Destino/pescado"
$ws.Range("B6").Value = "Madrid"
$ws.Range("C6").Value = "Barcelona"
$ws.Range("D6").Value = "Lisboa"

# Old quantities lived in B7:B9 - move them to D7:D9 (Lisboa column) and
# zero out Madrid / Barcelona.
$qty7 = $ws.Range("B7").Value()
$qty8 = $ws.Range("B8").Value()
$qty9 = $ws.Range("B9").Value()
$ws.Range("D7").Value = $qty7
$ws.Range("D8").Value = $qty8
$ws.Range("D9").Value = $qty9

$ws.Range("B7:C9").Value = 0

# Move the destination / distance helper table from D7:E9 to F7:G9.
$ws.Range("F6").Value = "Distancia"
$ws.Range("F7").Value = "Madrid"
$ws.Range("G7").Value = 800
$ws.Range("F8").Value = "Barcelona"
$ws.Range("G8").Value = 1100
$ws.Range("F9").Value = "Lisboa"
$ws.Range("G9").Value = 600

# The old D7:E9 "destination / distance" table content in column E is now
# empty since it moved to column G.
$ws.Range("E7:E9").ClearContents() | Out-Null

# --- Section A11:D17 - Ventas table: formulas now read the per-column
#     quantity (B/C/D) instead of being anchored to column B, and the
#     transport/depreciation formulas reference the relocated distance
#     table in column G. ---

$ws.Range("B12").Formula = "=B7*B2"
$ws.Range("C12:D12").Formula = "=C7*C2"

$ws.Range("B13:D13").Formula = "=B8*B3"

$ws.Range("B14:D14").Formula = "=B9*B4"

$ws.Range("B15").Formula = "=-IF(SUM(B12:B14)>0,5+2*`$G`$7,0)"
$ws.Range("C15").Formula = "=-IF(SUM(C12:C14)>0,5+2*`$G`$8,0)"
$ws.Range("D15").Formula = "=-IF(SUM(D12:D14)>0,5+2*`$G`$9,0)"

$ws.Range("B16").Formula = "=-SUM(B12:B14)*G7/100*0.01"
$ws.Range("C16").Formula = "=-SUM(C12:C14)*G8/100*0.01"
$ws.Range("D16").Formula = "=-SUM(D12:D14)*G9/100*0.01"

# New total-of-totals cell.
$ws.Range("E17").Formula = "=SUM(B17:D17)"

# The selection moved as part of editing the quantity table.
$ws.Range("D8").Select()

# Column A grows to fit the new, longer "Destino/pescado" label (the column
# already had Excel's best-fit autosize applied).
$ws.Columns("A:A").AutoFit() | Out-Null
